$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet tab text value occurrences (cell content), and the sheet name itself
$ws.Name = "Tienda - Velázquez"

# Insert a new column B ("serie") - this shifts old B:G to C:H
$ws.Columns.Item(2).Insert()

$ws.Range("B1").Value = "serie"

$ws.Range("A2:A11").Value = "Tienda - Velázquez"

$ws.Range("B2").Value = "V2"
$ws.Range("B3").Value = "V2"
$ws.Range("B4").Value = "V2"
$ws.Range("B5").Value = "V2"
$ws.Range("B6").Value = "V1"
$ws.Range("B7").Value = "V1"
$ws.Range("B8").Value = "V1"
$ws.Range("B9").Value = "V1"
$ws.Range("B10").Value = "V1"
$ws.Range("B11").Value = "V2"

$ws.Range("C11").Value = "BAR"
$ws.Range("E2").Value = "Mañana"
$ws.Range("E3").Value = "Mañana"
$ws.Range("E4").Value = "Mañana"
$ws.Range("E5").Value = "Mañana"
$ws.Range("E6").Value = "Mañana"
$ws.Range("E7").Value = "Mañana"
$ws.Range("E8").Value = "Mañana"
$ws.Range("E9").Value = "Mañana"
$ws.Range("E10").Value = "Mañana"
$ws.Range("E11").Value = "Mañana"

$ws.Range("F9").Value = "SMS"
$ws.Range("F10").Value = "TARJETA VISA"
$ws.Range("F11").Value = "TARJETA VISA"

$ws.Range("G2").Value = 72.75
$ws.Range("H2").Value = 9

$ws.Range("G3").Value = 202.1
$ws.Range("H3").Value = 24

$ws.Range("G4").Value = 83.40000000000001
$ws.Range("H4").Value = 9

$ws.Range("G5").Value = 139.8
$ws.Range("H5").Value = 14

$ws.Range("G6").Value = 806.3200000000001
$ws.Range("H6").Value = 77

$ws.Range("G7").Value = 2231.03
$ws.Range("H7").Value = 154

$ws.Range("G8").Value = 702.76
$ws.Range("H8").Value = 79

$ws.Range("G9").Value = 4.7
$ws.Range("H9").Value = 1

$ws.Range("G10").Value = 1809.8
$ws.Range("H10").Value = 154

$ws.Range("G11").Value = 93.2
$ws.Range("H11").Value = 13
